$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I19").Value = -0.4136749362997461
$ws.Range("J19").Value = 0.1266985367328312
$ws.Range("K19").Value = -0.1582366148292134
$ws.Range("L19").Value = 2.04042933805079
